$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.940.36"
$ws.Range("E2").Value = "  +2.69%  "

$ws.Range("D3").Value = "'3.064.65"
$ws.Range("E3").Value = "  +3.03%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'527.31"
$ws.Range("E5").Value = "  +6.47%  "

$ws.Range("D6").Value = "'143.23"
$ws.Range("E6").Value = "  +6.29%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  +5.86%  "

$ws.Range("D9").Value = "'7.67"
$ws.Range("E9").Value = "  +6.88%  "

$ws.Range("E10").Value = "  +7.83%  "

$ws.Range("E11").Value = "  +6.26%  "

$ws.Range("D13").Value = "'3.594.84"
$ws.Range("E13").Value = "  +3.20%  "

$ws.Range("D14").Value = "'27.47"
$ws.Range("E14").Value = "  +9.17%  "

$ws.Range("D15").Value = "'0.0000170"
$ws.Range("E15").Value = "  +16.79%  "

$ws.Range("D16").Value = "'57.951.06"
$ws.Range("E16").Value = "  +2.68%  "

$ws.Range("D17").Value = "'6.22"
$ws.Range("E17").Value = "  +7.31%  "

$ws.Range("D18").Value = "'3.068.79"
$ws.Range("E18").Value = "  +3.00%  "

$ws.Range("D19").Value = "'13.22"
$ws.Range("E19").Value = "  +7.53%  "

$ws.Range("D20").Value = "'8.19"
$ws.Range("E20").Value = "  +5.61%  "

$ws.Range("D21").Value = "'341.45"
$ws.Range("E21").Value = "  +5.12%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "'5.68"
$ws.Range("E23").Value = "  -0.82%  "

$ws.Range("D24").Value = "'0.503"
$ws.Range("E24").Value = "  +7.37%  "

$ws.Range("D25").Value = "'64.95"
$ws.Range("E25").Value = "  +5.68%  "

$ws.Range("D26").Value = "'0.171"
$ws.Range("E26").Value = "  +6.68%  "

$ws.Range("D27").Value = "'0.0₃0979"
$ws.Range("E27").Value = "  +10.19%  "

$ws.Range("E28").Value = "  +0.47%  "

$ws.Range("D29").Value = "'6.96"
$ws.Range("E29").Value = "  +7.20%  "

$ws.Range("D30").Value = "'7.42"
$ws.Range("E30").Value = "  +10.19%  "

$ws.Range("E31").Value = "  +6.92%  "

$ws.Range("E32").Value = "  +7.76%  "

$ws.Range("D33").Value = "'21.10"
$ws.Range("E33").Value = "  +4.00%  "

$ws.Range("D34").Value = "'4.79"
$ws.Range("E34").Value = "  +7.75%  "

$ws.Range("D35").Value = "'157.26"
$ws.Range("E35").Value = "  +3.48%  "

$ws.Range("D36").Value = "'6.01"
$ws.Range("E36").Value = "  +7.81%  "

$ws.Range("D37").Value = "'1.33"
$ws.Range("E37").Value = "  +4.67%  "

$ws.Range("D38").Value = "'26.49"
$ws.Range("E38").Value = "  +15.31%  "

$ws.Range("D39").Value = "'0.0707"
$ws.Range("E39").Value = "  +5.84%  "

$ws.Range("D40").Value = "'3.101.73"
$ws.Range("E40").Value = "  +3.20%  "

$ws.Range("D41").Value = "'37.87"
$ws.Range("E41").Value = "  +3.34%  "

$ws.Range("D42").Value = "'3.91"
$ws.Range("E42").Value = "  +10.26%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("E44").Value = "  +4.77%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.48"
$ws.Range("E45").Value = "  +6.57%  "

$ws.Range("D46").Value = "'2.337.69"
$ws.Range("E46").Value = "  +5.19%  "

$ws.Range("D47").Value = "'1.03"
$ws.Range("E47").Value = "  +3.69%  "

$ws.Range("D48").Value = "'2.01"
$ws.Range("E48").Value = "  +3.93%  "

$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'6.05"
$ws.Range("E49").Value = "  +5.68%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0245"
$ws.Range("E50").Value = "  +3.73%  "

$ws.Range("D51").Value = "'20.22"
$ws.Range("E51").Value = "  +7.25%  "
